$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 117
$ws.Range("A117").Value = 7
$ws.Range("B117").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C117").Value = "Ñuble"
$ws.Range("D117").Value = 44418
$ws.Range("D117").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E117").Value = 16
$ws.Range("F117").Value = "Fruta"
$ws.Range("G117").Value = 100108
$ws.Range("H117").Value = "Tropicales y subtropicales"
$ws.Range("I117").Value = 100108005
$ws.Range("J117").Value = "Piña"
$ws.Range("K117").Value = "Caramelo"
$ws.Range("L117").Value = "Primera"
$ws.Range("M117").Value = 100
$ws.Range("N117").Value = 18000
$ws.Range("O117").Value = 19000
$ws.Range("P117").Value = 18500
$ws.Range("Q117").Value = "$/caja 12 unidades"
$ws.Range("R117").Value = "Ecuador"
$ws.Range("S117").Value = 1542
$ws.Range("T117").Value = 12

# Row 118
$ws.Range("A118").Value = 7
$ws.Range("B118").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C118").Value = "Ñuble"
$ws.Range("D118").Value = 44418
$ws.Range("D118").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E118").Value = 16
$ws.Range("F118").Value = "Fruta"
$ws.Range("G118").Value = 100108
$ws.Range("H118").Value = "Tropicales y subtropicales"
$ws.Range("I118").Value = 100108005
$ws.Range("J118").Value = "Piña"
$ws.Range("K118").Value = "Caramelo"
$ws.Range("L118").Value = "Segunda"
$ws.Range("M118").Value = 80
$ws.Range("N118").Value = 18000
$ws.Range("O118").Value = 19000
$ws.Range("P118").Value = 18500
$ws.Range("Q118").Value = "$/caja 14 unidades"
$ws.Range("R118").Value = "Ecuador"
$ws.Range("S118").Value = 1321
$ws.Range("T118").Value = 14
